$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Silver Rear_side")
$ws1.Range("B20").NumberFormat = "@"
$ws1.Range("B20").Value = "5,399"

$ws2 = $wb.Worksheets.Item("Silver Busbar front-side")
$ws2.Range("B20").NumberFormat = "@"
$ws2.Range("B20").Value = "8,083"

$ws3 = $wb.Worksheets.Item("Silver finger front-side")
$ws3.Range("B20").NumberFormat = "@"
$ws3.Range("B20").Value = "8,133"

$ws4 = $wb.Worksheets.Item("USD_CNY")
$ws4.Range("B20").NumberFormat = "@"
$ws4.Range("B20").Value = "7.2717"
